$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-QuotedText($row, $text) {
    $ws.Range("Z1").Formula = "=CHAR(39)&""$text""&CHAR(39)"
    $ws.Range("Z1").Copy()
    $ws.Range("B$row").PasteSpecial(-4163)
}

Set-QuotedText 2 "Akkermansia_muciniphila_ATCC_BAA_835.mat"
$ws.Range("C2").Value = 0

Set-QuotedText 3 "Alistipes_finegoldii_DSM_17242.mat"
$ws.Range("C3").Value = 0

Set-QuotedText 4 "Alistipes_indistinctus_YIT_12060.mat"
$ws.Range("C4").Value = 0

Set-QuotedText 5 "Alistipes_putredinis_DSM_17216.mat"
$ws.Range("C5").Value = 0

Set-QuotedText 6 "Alistipes_shahii_WAL_8301.mat"
$ws.Range("C6").Value = 0

Set-QuotedText 7 "Bacteroides_cellulosilyticus_DSM_14838.mat"
$ws.Range("C7").Value = 0

Set-QuotedText 8 "Bacteroides_fragilis_3_1_12.mat"
$ws.Range("C8").Value = 0

Set-QuotedText 9 "Bacteroides_oleiciplenus_YIT_12058.mat"
$ws.Range("C9").Value = 0

Set-QuotedText 10 "Bacteroides_ovatus_ATCC_8483.mat"
$ws.Range("C10").Value = 0

Set-QuotedText 11 "Bacteroides_plebeius_M12_DSM_17135.mat"
$ws.Range("C11").Value = 0

Set-QuotedText 12 "Bacteroides_salyersiae_WAL_10018.mat"
$ws.Range("C12").Value = 0

Set-QuotedText 13 "Bacteroides_thetaiotaomicron_VPI_5482.mat"
$ws.Range("C13").Value = 0

Set-QuotedText 14 "Bacteroides_uniformis_ATCC_8492.mat"
$ws.Range("C14").Value = 0

Set-QuotedText 15 "Bacteroides_vulgatus_ATCC_8482.mat"
$ws.Range("C15").Value = 0

Set-QuotedText 16 "Barnesiella_intestinihominis_YIT_11860.mat"
$ws.Range("C16").Value = 0

Set-QuotedText 17 "Bifidobacterium_animalis_lactis_AD011.mat"
$ws.Range("C17").Value = 0

Set-QuotedText 18 "Bilophila_wadsworthia_3_1_6.mat"
$ws.Range("C18").Value = 0

Set-QuotedText 19 "Escherichia_coli_O157_H7_str_Sakai.mat"
$ws.Range("C19").Value = 0

Set-QuotedText 20 "Eubacterium_limosum_KIST612.mat"
$ws.Range("C20").Value = 0

Set-QuotedText 21 "Eubacterium_ramulus_ATCC_29099.mat"
$ws.Range("C21").Value = 0

Set-QuotedText 22 "Flavonifractor_plautii_ATCC_29863.mat"
$ws.Range("C22").Value = 0.004

Set-QuotedText 23 "Marvinbryantia_formatexigens_I_52_DSM_14469.mat"
$ws.Range("C23").Value = 0

Set-QuotedText 24 "Odoribacter_splanchnicus_1651_6_DSM_20712.mat"
$ws.Range("C24").Value = 0

Set-QuotedText 25 "Parabacteroides_distasonis_ATCC_8503.mat"
$ws.Range("C25").Value = 0

Set-QuotedText 26 "Parabacteroides_johnsonii_DSM_18315.mat"
$ws.Range("C26").Value = 0

Set-QuotedText 27 "Paraprevotella_xylaniphila_YIT_11841.mat"
$ws.Range("C27").Value = 0

Set-QuotedText 28 "Parasutterella_excrementihominis_YIT_11859.mat"
$ws.Range("C28").Value = 0.032

Set-QuotedText 29 "Phascolarctobacterium_succinatutens_YIT_12067.mat"
$ws.Range("C29").Value = 0.075

Set-QuotedText 30 "Prevotella_copri_CB7_DSM_18205.mat"
$ws.Range("C30").Value = -0

Set-QuotedText 31 "Prevotella_stercorea_DSM_18206.mat"
$ws.Range("C31").Value = 0.024

Set-QuotedText 32 "Roseburia_inulinivorans_DSM_16841.mat"
$ws.Range("C32").Value = 0.133

Set-QuotedText 33 "Sutterella_wadsworthensis_3_1_45B.mat"
$ws.Range("C33").Value = 0.732

$ws.Range("Z1").Value = ""
$ws.Application.CutCopyMode = $false
